{"js": "// Update the date heading and the 25 \"three-digit \u00d7 one-digit\" problems\n// in the practice table. Cells are addressed by (row, column) position\n// rather than by searching for the old text, because several problems\n// share the same text (e.g. \"969\u00d74=\" appears twice with different\n// replacements), so a plain text search/replace would be ambiguous.\n\nconst body = context.document.body;\n\n// 1. Update the date/weekday heading (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.insertText(\"2024-06-10 Monday\", \"Replace\");\n\n// 2. Update the math problems inside the table.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index -> new cell values (left to right). Only rows that hold\n// problems carry text; the other rows are blank work-space rows.\nconst rowUpdates = {\n  0: [\"767\u00d75=\", \"959\u00d75=\", \"601\u00d79=\", \"336\u00d76=\", \"596\u00d73=\"],\n  4: [\"979\u00d78=\", \"583\u00d72=\", \"130\u00d74=\", \"527\u00d77=\", \"781\u00d77=\"],\n  9: [\"646\u00d72=\", \"455\u00d74=\", \"466\u00d79=\", \"978\u00d72=\", \"984\u00d77=\"],\n  14: [\"279\u00d79=\", \"466\u00d77=\", \"867\u00d78=\", \"220\u00d72=\", \"316\u00d77=\"],\n  19: [\"613\u00d77=\", \"319\u00d73=\", \"811\u00d73=\", \"617\u00d79=\", \"267\u00d72=\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const values = rowUpdates[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(Number(rowIndex), col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 \"three-digit x one-digit\" problems in\n# the practice table. Cells are addressed by (row, column) position rather\n# than by searching for the old text, because several problems share the\n# same text (e.g. \"969x4=\" appears twice with different replacements), so a\n# plain Find/Replace would be ambiguous.\n\n$d = $word.ActiveDocument\n\n# 1. Update the date/weekday heading (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-10 Monday\"\n\n# 2. Update the math problems inside the table (1-based row/column index,\n# counting every table row including the blank work-space rows).\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1, 1).Range.Text = \"767\u00d75=\"\n$tbl.Cell(1, 2).Range.Text = \"959\u00d75=\"\n$tbl.Cell(1, 3).Range.Text = \"601\u00d79=\"\n$tbl.Cell(1, 4).Range.Text = \"336\u00d76=\"\n$tbl.Cell(1, 5).Range.Text = \"596\u00d73=\"\n\n$tbl.Cell(5, 1).Range.Text = \"979\u00d78=\"\n$tbl.Cell(5, 2).Range.Text = \"583\u00d72=\"\n$tbl.Cell(5, 3).Range.Text = \"130\u00d74=\"\n$tbl.Cell(5, 4).Range.Text = \"527\u00d77=\"\n$tbl.Cell(5, 5).Range.Text = \"781\u00d77=\"\n\n$tbl.Cell(10, 1).Range.Text = \"646\u00d72=\"\n$tbl.Cell(10, 2).Range.Text = \"455\u00d74=\"\n$tbl.Cell(10, 3).Range.Text = \"466\u00d79=\"\n$tbl.Cell(10, 4).Range.Text = \"978\u00d72=\"\n$tbl.Cell(10, 5).Range.Text = \"984\u00d77=\"\n\n$tbl.Cell(15, 1).Range.Text = \"279\u00d79=\"\n$tbl.Cell(15, 2).Range.Text = \"466\u00d77=\"\n$tbl.Cell(15, 3).Range.Text = \"867\u00d78=\"\n$tbl.Cell(15, 4).Range.Text = \"220\u00d72=\"\n$tbl.Cell(15, 5).Range.Text = \"316\u00d77=\"\n\n$tbl.Cell(20, 1).Range.Text = \"613\u00d77=\"\n$tbl.Cell(20, 2).Range.Text = \"319\u00d73=\"\n$tbl.Cell(20, 3).Range.Text = \"811\u00d73=\"\n$tbl.Cell(20, 4).Range.Text = \"617\u00d79=\"\n$tbl.Cell(20, 5).Range.Text = \"267\u00d72=\"\n"}
